$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Swap the contents of columns C (codeforiati:group-code) and D (codeforiati:group-name)
# for every row, including the header row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
